$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing data row (row 46) into the new row 47 so the
# new row inherits the exact same cell style used by every other data row
# in the sheet (General number format, center/center alignment, no quote
# prefix).
$ws.Range("A46:C46").Copy($ws.Range("A47:C47"))

# C47 is about to be overwritten with a number anyway, so use it as scratch
# space: typing a leading apostrophe forces Excel to store "2025/12/27" as
# plain text instead of re-interpreting it as a date. A Paste Special that
# copies values only (not formats) then transfers just that text content -
# without the quote-prefix formatting - onto the already correctly styled
# A47 cell.
$ws.Cells.Item(47, 3).Formula = "'2025/12/27"
$ws.Cells.Item(47, 3).Copy()
$ws.Cells.Item(47, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

# Fill in the remaining real values for the new row.
$ws.Cells.Item(47, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(47, 3).Value = 1104
